# Generate Report for Handback
# - Renames the in-flight handback test file (759239cf-...) to a freshly
#   generated one (4fda43a1-...) across all three sheets, refreshing its
#   timestamps/hashes.
# - Appends a brand-new second handback row (e0d321a8-...) to all three
#   sheets / tables.

$wb = $excel.ActiveWorkbook

$oldBase = "759239cf-2818-4f34-9c84-0fae4df38b1c"
$newBase1 = "4fda43a1-c97f-474b-9818-382765714de1"   # renamed row (was row 2)
$newBase2 = "e0d321a8-fdee-4cdc-9c96-bb151ac53a11"   # brand-new row (row 3)

$hashZh1 = "6b8d52bc3f1a85ba976901c091ae6e4e4ee85974"
$hashZh2 = "4ab91fcfe92d42d533f96599be05274f07e1e08e"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null

# Row 2 (renamed file)
$wsO.Range("A2").Value = "$newBase1.md"
$wsO.Range("B2").Value = "e2e\$newBase1.md"
$wsO.Range("C2").Value = ".md"
$wsO.Range("E2").Value = $statusText
$wsO.Range("F2").Value = $statusText
$wsO.Range("G2").Value = "2016-08-19 15:07:47"

$wsO.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 (new file)
$wsO.Range("A3").Value = "$newBase2.md"
$wsO.Range("B3").Value = "e2e\$newBase2.md"
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = $statusText
$wsO.Range("F3").Value = $statusText
$wsO.Range("G3").Value = "2016-08-19 15:07:47"
$wsO.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsO.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase1.md", "", "", "e2e\$newBase1.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase2.md", "", "", "e2e\$newBase2.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")

$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add() | Out-Null

# Row 2 (renamed file)
$wsZ.Range("A2").Value = "$newBase1.md"
$wsZ.Range("B2").Value = ".md"
$wsZ.Range("C2").Value = $statusText
$wsZ.Range("D2").Value = "e2e"
$wsZ.Range("E2").Value = "ht"
$wsZ.Range("F2").Value = "'False"
$wsZ.Range("G2").Value = "$newBase1.$hashZh1.zh-cn.xlf"
$wsZ.Range("H2").Value = "2016-08-19 15:07:43"
$wsZ.Range("I2").Value = "$newBase1.md"
$wsZ.Range("J2").Value = "$newBase1.$hashZh1.zh-cn.xlf"
$wsZ.Range("K2").Value = "2016-08-19 15:08:03"
$wsZ.Range("L2").Value = "'"
$wsZ.Range("M2").Value = "'True"
$wsZ.Range("N2").Value = "'"
$wsZ.Range("O2").Value = "'False"
$wsZ.Range("P2").Value = "'"
$wsZ.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 (new file)
$wsZ.Range("A3").Value = "$newBase2.md"
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = $statusText
$wsZ.Range("D3").Value = "e2e"
$wsZ.Range("E3").Value = "ht"
$wsZ.Range("F3").Value = "'True"
$wsZ.Range("G3").Value = "$newBase2.$hashZh2.zh-cn.xlf"
$wsZ.Range("H3").Value = "2016-08-19 15:07:43"
$wsZ.Range("I3").Value = "$newBase2.md"
$wsZ.Range("J3").Value = "$newBase2.$hashZh2.zh-cn.xlf"
$wsZ.Range("K3").Value = "2016-08-19 15:08:03"
$wsZ.Range("L3").Value = "'"
$wsZ.Range("M3").Value = "'True"
$wsZ.Range("N3").Value = "'"
$wsZ.Range("O3").Value = "'False"
$wsZ.Range("P3").Value = "'"
$wsZ.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZ.Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f5d5f5d0f92d2f2bdf7b38d183089e46b2225f04/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f5d5f5d0f92d2f2bdf7b38d183089e46b2225f04/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")

$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add() | Out-Null

# Row 2 (renamed file)
$wsD.Range("A2").Value = "$newBase1.md"
$wsD.Range("B2").Value = ".md"
$wsD.Range("C2").Value = $statusText
$wsD.Range("D2").Value = "e2e"
$wsD.Range("E2").Value = "ht"
$wsD.Range("F2").Value = "'False"
$wsD.Range("G2").Value = "$newBase1.$hashZh1.de-de.xlf"
$wsD.Range("H2").Value = "2016-08-19 15:07:47"
$wsD.Range("I2").Value = "$newBase1.md"
$wsD.Range("J2").Value = "$newBase1.$hashZh1.de-de.xlf"
$wsD.Range("K2").Value = "2016-08-19 15:08:16"
$wsD.Range("L2").Value = "'"
$wsD.Range("M2").Value = "'True"
$wsD.Range("N2").Value = "'"
$wsD.Range("O2").Value = "'False"
$wsD.Range("P2").Value = "'"
$wsD.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 (new file)
$wsD.Range("A3").Value = "$newBase2.md"
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = $statusText
$wsD.Range("D3").Value = "e2e"
$wsD.Range("E3").Value = "ht"
$wsD.Range("F3").Value = "'True"
$wsD.Range("G3").Value = "$newBase2.$hashZh2.de-de.xlf"
$wsD.Range("H3").Value = "2016-08-19 15:07:47"
$wsD.Range("I3").Value = "$newBase2.md"
$wsD.Range("J3").Value = "$newBase2.$hashZh2.de-de.xlf"
$wsD.Range("K3").Value = "2016-08-19 15:08:16"
$wsD.Range("L3").Value = "'"
$wsD.Range("M3").Value = "'True"
$wsD.Range("N3").Value = "'"
$wsD.Range("O3").Value = "'False"
$wsD.Range("P3").Value = "'"
$wsD.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsD.Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1be1a1a7607c57a55d48c40d1d64b0760c18ce55/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1be1a1a7607c57a55d48c40d1d64b0760c18ce55/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null

Write-Host "Handback report rows updated."
